$wb = $excel.ActiveWorkbook

# ColumnWidth values are quantized to whole pixels by this runtime, so the
# nearest reachable width to the target 17.2159881591797 (characters) is
# produced by any ColumnWidth input in [16.2501, 16.4166); 16.3 is used.
$targetColumnWidth = 16.3

# --- Overview sheet: status + timestamp refreshed for new handoff ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 19:08:23"
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# --- zh-cn sheet: status + handoff timestamp refreshed ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 19:08:18"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# --- de-de sheet: status + (shared) handoff timestamp refreshed ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 19:08:23"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
